# fcr collection current_config_extract separated into sections
#
# This script reproduces, via the Excel COM object model, the changes
# captured in the target OOXML diff:
#  - report_requisites!A1 ("Contents" header / hyperlink) is extended
#    across the whole header row (B1:BC1), matching A1's style and
#    hyperlink target.
#  - report_requisites!N3 ("value" label) moves to F3.
#  - project_steps: the flag cell (value 1) that used to sit in column H
#    of row 25 moves to column H of row 43 (a row that used to have no
#    value there) — i.e. the "current_config_extract" row moved to a
#    later/ different section.
#  - Various saved cursor/selection positions + scroll/frozen-pane state
#    change on several sheets, and the active worksheet tab changes from
#    report_requisites to project_steps.

$wb = $excel.ActiveWorkbook

$wsContents   = $wb.Worksheets.Item("contents")
$wsRequisites = $wb.Worksheets.Item("report_requisites")
$wsSteps      = $wb.Worksheets.Item("project_steps")
$wsLinks      = $wb.Worksheets.Item("io_data_names_links")
$wsNames      = $wb.Worksheets.Item("in_out_data_names")

# ---------------------------------------------------------------------
# contents: move saved selection from D58 to A3 (sheet stays inactive)
# ---------------------------------------------------------------------
$wsContents.Range("A3").Select()

# ---------------------------------------------------------------------
# report_requisites: spread the "Contents" hyperlink header across
# B1:BC1 (same text/style/hyperlink as A1), and move the "value" label
# cell from N3 to F3.
# ---------------------------------------------------------------------
$headerRange = $wsRequisites.Range("B1:BC1")
$wsRequisites.Hyperlinks.Add($headerRange, "", "contents!A1", "", "Contents") | Out-Null
$headerRange.Style = "Hyperlink"   # re-apply A1's exact (shared) cell style
$headerRange.Value() = $wsRequisites.Range("A1").Value()

$wsRequisites.Range("A3").Copy()
$wsRequisites.Range("F3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats - reuses A3's exact style
$excel.CutCopyMode = $false
$wsRequisites.Range("F3").Value() = $wsRequisites.Range("N3").Value()
$wsRequisites.Range("N3").Clear()

# Saved cursor position for report_requisites moves too.
$wsRequisites.Range("E34").Select()

# ---------------------------------------------------------------------
# project_steps: the "current_config_extract" flag (value 1) moves from
# H25 to H43 - the section it documents was split out separately.
# ---------------------------------------------------------------------
$wsSteps.Range("H25").Clear()
$wsSteps.Range("H43").Value = 1

# ---------------------------------------------------------------------
# io_data_names_links / in_out_data_names: saved cursor positions move.
# ---------------------------------------------------------------------
$wsLinks.Range("B14").Select()
$wsNames.Range("B18").Select()

# ---------------------------------------------------------------------
# project_steps becomes the active sheet/tab, with a new saved
# selection (F43:H43) inside its frozen pane.
# ---------------------------------------------------------------------
$wsSteps.Activate()
$wsSteps.Range("F43:H43").Select()
